$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.154.17'
$ws.Range("E2").Value = '  -0.86%  '
$ws.Range("D3").Value = '3.273.32'
$ws.Range("E3").Value = '  +0.42%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.89'
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.80'
$ws.Range("E6").Value = '  +1.71%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +1.63%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.130'
$ws.Range("E9").Value = '  -1.71%  '
$ws.Range("E10").Value = '  -1.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.410'
$ws.Range("E11").Value = '  -2.62%  '
$ws.Range("D12").Value = '3.841.26'
$ws.Range("E12").Value = '  +0.18%  '
$ws.Range("E13").Value = '  +0.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.50'
$ws.Range("E14").Value = '  -3.18%  '
$ws.Range("D15").Value = '68.147.39'
$ws.Range("E15").Value = '  -0.80%  '
$ws.Range("E16").Value = '  -1.10%  '
$ws.Range("D17").Value = '3.275.19'
$ws.Range("E17").Value = '  +0.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.72'
$ws.Range("E18").Value = '  -1.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.44'
$ws.Range("E19").Value = '  -0.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '413.21'
$ws.Range("E20").Value = '  +5.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.56'
$ws.Range("E21").Value = '  -1.16%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.45'
$ws.Range("E23").Value = '  -0.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.510'
$ws.Range("E24").Value = '  -0.88%  '
$ws.Range("E25").Value = '  -1.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.187'
$ws.Range("E26").Value = '  -2.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.48'
$ws.Range("E27").Value = '  -0.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").Value = '  -0.26%  '
$ws.Range("E29").Value = '  -1.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.72'
$ws.Range("E30").Value = '  -0.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.49'
$ws.Range("E31").Value = '  -3.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.90'
$ws.Range("E32").Value = '  -2.79%  '
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("E34").Value = '  -2.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '164.58'
$ws.Range("E35").Value = '  +0.72%  '
$ws.Range("E36").Value = '  -2.35%  '
$ws.Range("E37").Value = '  -1.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '27.31'
$ws.Range("E38").Value = '  +4.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.800'
$ws.Range("E39").Value = '  -3.24%  '
$ws.Range("E40").Value = '  -1.91%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.36'
$ws.Range("E41").Value = '  -3.60%  '
$ws.Range("D42").Value = '2.669.34'
$ws.Range("E42").Value = '  +2.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.83'
$ws.Range("E43").Value = '  -1.17%  '
$ws.Range("B44").Value = 'Hedera'
$ws.Range("C44").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0680'
$ws.Range("E44").Value = '  -1.07%  '
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.44'
$ws.Range("E45").Value = '  -1.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '337.82'
$ws.Range("E46").Value = '  -1.81%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.58'
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0274'
$ws.Range("E48").Value = '  -2.62%  '
$ws.Range("E49").Value = '  +0.45%  '
$ws.Range("E50").Value = '  -0.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.973'
$ws.Range("E51").Value = '  -0.58%  '
